{"js": "const pairs = [\n  [\"772\u00f79=\", \"791\u00f73=\"],\n  [\"676\u00f79=\", \"742\u00f76=\"],\n  [\"573\u00f78=\", \"367\u00f75=\"],\n  [\"622\u00f74=\", \"536\u00f76=\"],\n  [\"442\u00f76=\", \"109\u00f78=\"],\n  [\"949\u00f74=\", \"211\u00f76=\"],\n  [\"449\u00f72=\", \"376\u00f78=\"],\n  [\"161\u00f73=\", \"907\u00f76=\"],\n  [\"936\u00f75=\", \"850\u00f72=\"],\n  [\"628\u00f77=\", \"447\u00f74=\"],\n  [\"296\u00f75=\", \"297\u00f74=\"],\n  [\"380\u00f73=\", \"360\u00f76=\"],\n  [\"674\u00f72=\", \"554\u00f77=\"],\n  [\"178\u00f78=\", \"387\u00f77=\"],\n  [\"876\u00f77=\", \"342\u00f77=\"],\n  [\"389\u00f73=\", \"161\u00f76=\"],\n  [\"778\u00f72=\", \"565\u00f75=\"],\n  [\"955\u00f73=\", \"220\u00f73=\"],\n  [\"280\u00f72=\", \"139\u00f79=\"],\n  [\"657\u00f72=\", \"434\u00f74=\"],\n  [\"387\u00f78=\", \"895\u00f76=\"],\n  [\"448\u00f72=\", \"508\u00f79=\"],\n  [\"225\u00f78=\", \"165\u00f79=\"],\n  [\"534\u00f75=\", \"365\u00f78=\"],\n  [\"103\u00f79=\", \"270\u00f72=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$pairs = @(\n  @(\"772\u00f79=\", \"791\u00f73=\"),\n  @(\"676\u00f79=\", \"742\u00f76=\"),\n  @(\"573\u00f78=\", \"367\u00f75=\"),\n  @(\"622\u00f74=\", \"536\u00f76=\"),\n  @(\"442\u00f76=\", \"109\u00f78=\"),\n  @(\"949\u00f74=\", \"211\u00f76=\"),\n  @(\"449\u00f72=\", \"376\u00f78=\"),\n  @(\"161\u00f73=\", \"907\u00f76=\"),\n  @(\"936\u00f75=\", \"850\u00f72=\"),\n  @(\"628\u00f77=\", \"447\u00f74=\"),\n  @(\"296\u00f75=\", \"297\u00f74=\"),\n  @(\"380\u00f73=\", \"360\u00f76=\"),\n  @(\"674\u00f72=\", \"554\u00f77=\"),\n  @(\"178\u00f78=\", \"387\u00f77=\"),\n  @(\"876\u00f77=\", \"342\u00f77=\"),\n  @(\"389\u00f73=\", \"161\u00f76=\"),\n  @(\"778\u00f72=\", \"565\u00f75=\"),\n  @(\"955\u00f73=\", \"220\u00f73=\"),\n  @(\"280\u00f72=\", \"139\u00f79=\"),\n  @(\"657\u00f72=\", \"434\u00f74=\"),\n  @(\"387\u00f78=\", \"895\u00f76=\"),\n  @(\"448\u00f72=\", \"508\u00f79=\"),\n  @(\"225\u00f78=\", \"165\u00f79=\"),\n  @(\"534\u00f75=\", \"365\u00f78=\"),\n  @(\"103\u00f79=\", \"270\u00f72=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
